$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the displayed text of the hyperlink cells in column D.
# The underlying hyperlink relationships are left untouched; only the
# shared-string text shown in each cell changes.
$ws.Range("D9").Value  = "https://www.chadjemmett.com1"
$ws.Range("D10").Value = "https://www.chadjemmett.com2"
$ws.Range("D11").Value = "https://www.chadjemmett.com3"
$ws.Range("D12").Value = "https://www.chadjemmett.com4"
$ws.Range("D13").Value = "https://www.chadjemmett.com4"
$ws.Range("D14").Value = "https://www.chadjemmett.com5"

$ws.Range("D18").Value = "www.google.com"
$ws.Range("D19").Value = "www.boinggoin.com"
$ws.Range("D20").Value = "www.facebooko.com"
$ws.Range("D21").Value = "www.metafileter.com"
$ws.Range("D22").Value = "www.bobo.com"
$ws.Range("D23").Value = "www.twitch.com"

# Reflect the new selected cell recorded in the saved workbook view.
$ws.Range("D24").Select()
